$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with revised totals
$ws.Range("C253").Value = 93020.38
$ws.Range("C255").Value = 102609.33

# Row 256: change id_loja, total_venda and nome (store)
$ws.Range("B256").Value = 1
$ws.Range("C256").Value = 17144.8
$ws.Range("E256").Value = "Bibi Cell Mundi"

# Row 257: change id_loja, total_venda and nome (store)
$ws.Range("B257").Value = 2
$ws.Range("C257").Value = 4714
$ws.Range("E257").Value = "Bibi Cell Manauara"

# New row 258
$ws.Range("A258").Value = 6
$ws.Range("B258").Value = 3
$ws.Range("C258").Value = 4464
$ws.Range("D258").Value = 2025
$ws.Range("E258").Value = "Bibi Cell Vieiralves"

# New row 259
$ws.Range("A259").Value = 6
$ws.Range("B259").Value = 4
$ws.Range("C259").Value = 6470.01
$ws.Range("D259").Value = 2025
$ws.Range("E259").Value = "Bibi Cell Ponta Negra"
